$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.134.36'
$ws.Range("E2").Value = '  -3.98%  '
$ws.Range("D3").Value = '2.443.96'
$ws.Range("E3").Value = '  -4.29%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = "'525.63"
$ws.Range("E5").Value = '  -3.35%  '
$ws.Range("D6").Value = "'133.40"
$ws.Range("E6").Value = '  -9.30%  '
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").Value = "'0.551"
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("D9").Value = '2.449.87'
$ws.Range("E9").Value = '  -5.44%  '
$ws.Range("D10").Value = "'0.0988"
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = "'5.30"
$ws.Range("E12").Value = '  -4.40%  '
$ws.Range("D13").Value = "'0.341"
$ws.Range("E13").Value = '  -6.82%  '
$ws.Range("D14").Value = '2.884.24'
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").Value = '58.013.63'
$ws.Range("E15").Value = '  -4.14%  '
$ws.Range("D16").Value = "'22.34"
$ws.Range("E16").Value = '  -9.14%  '
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = '  -5.05%  '
$ws.Range("D18").Value = '2.455.28'
$ws.Range("E18").Value = '  -4.30%  '
$ws.Range("D19").Value = "'10.58"
$ws.Range("E19").Value = '  -7.39%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'317.65"
$ws.Range("E20").Value = '  -3.76%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = "'4.14"
$ws.Range("E21").Value = '  -5.69%  '
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = "'5.64"
$ws.Range("E23").Value = '  -5.93%  '
$ws.Range("D24").Value = "'61.84"
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("D25").Value = "'0.402"
$ws.Range("E25").Value = '  -9.70%  '
$ws.Range("E26").Value = '  -3.21%  '
$ws.Range("D27").Value = "'0.984"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "'7.45"
$ws.Range("E28").Value = '  -8.08%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0748'
$ws.Range("E29").Value = '  -8.42%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = "'6.47"
$ws.Range("E30").Value = '  -9.93%  '
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = '  -5.07%  '
$ws.Range("D32").Value = "'162.96"
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").Value = "'1.05"
$ws.Range("E34").Value = '  -14.05%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'18.08"
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'1.34"
$ws.Range("E36").Value = '  -9.74%  '
$ws.Range("E37").Value = '  -11.86%  '
$ws.Range("D38").Value = "'1.52"
$ws.Range("E38").Value = '  -8.51%  '
$ws.Range("D39").Value = "'36.31"
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("E40").Value = '  -7.32%  '
$ws.Range("D41").Value = "'0.775"
$ws.Range("E41").Value = '  -9.08%  '
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").Value = "'269.31"
$ws.Range("E43").Value = '  -13.18%  '
$ws.Range("D44").Value = "'4.96"
$ws.Range("E44").Value = '  -14.11%  '
$ws.Range("D45").Value = "'10.82"
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").Value = "'0.581"
$ws.Range("E46").Value = '  -5.39%  '
$ws.Range("D47").Value = "'0.0921"
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("D48").Value = "'119.20"
$ws.Range("E48").Value = '  -6.12%  '
$ws.Range("E49").Value = '  -5.47%  '
$ws.Range("D50").Value = "'0.0216"
$ws.Range("E50").Value = '  -7.05%  '
$ws.Range("D51").Value = "'16.84"
$ws.Range("E51").Value = '  -9.54%  '
